# Regenerate the "K" (strikeout) column (column G) values for rows 2-79
# using the freshly computed s_vals (derived from the underlying play-by-play
# data rather than the previous Strike# based count).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G ("K") values for rows 2 through 79, in row order.
$kValues = @(
    3,3,0,1,1,1,2,0,0,1,
    0,2,0,0,2,2,2,2,1,1,
    0,0,2,1,0,2,1,3,0,1,
    1,1,1,0,1,1,1,1,2,1,
    1,0,1,1,1,1,1,1,0,0,
    1,1,1,0,1,0,1,0,0,1,
    0,1,0,0,1,2,1,1,1,1,
    2,2,1,1,2,3,0,1
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
